$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's (2022-09-20) trading row for DELHIVERY at row 84,
# following the same layout as the existing historical rows.
$ws.Range("A84").Value = 44824
$ws.Range("A84").NumberFormat = "YYYY-MM-DD"
$ws.Range("B84").Value = "DELHIVERY"
$ws.Range("C84").Value = "EQ"
$ws.Range("D84").Value = 575.55
$ws.Range("E84").Value = 588.9
$ws.Range("F84").Value = 588.9
$ws.Range("G84").Value = 576.3
$ws.Range("H84").Value = 582
$ws.Range("I84").Value = 581.35
$ws.Range("J84").Value = 581.24
$ws.Range("K84").Value = 292812
$ws.Range("L84").Value = 17019475625000
$ws.Range("M84").Value = 10777
$ws.Range("N84").Value = 142052
$ws.Range("O84").Value = 0.4851
